$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (price/volume updates)
$updates = @{
    'D2' = '29.272.02'
    'E2' = '  +0.27%  '
    'D3' = '1.901.74'
    'E3' = '  -0.46%  '
    'D4' = '1.002'
    'E4' = '  -0.18%  '
    'D5' = '326.73'
    'E5' = '  -0.38%  '
    'E6' = '  -0.18%  '
    'D7' = '0.4654'
    'E7' = '  -0.10%  '
    'D8' = '0.3927'
    'E8' = '  -0.07%  '
    'D9' = '0.07888'
    'E9' = '  -1.13%  '
    'D10' = '0.9889'
    'E10' = '  -1.59%  '
    'D11' = '21.99'
    'E11' = '  -1.69%  '
    'D12' = '1.923.89'
    'E12' = '  -2.29%  '
    'D13' = '7.079'
    'E13' = '  -1.01%  '
    'D14' = '5.754'
    'E14' = '  -0.76%  '
    'D15' = '0.06998'
    'E15' = '  +0.38%  '
    'D16' = '88.40'
    'E16' = '  -0.47%  '
    'D17' = '1.002'
    'E17' = '  -0.21%  '
    'D18' = '0.000009992'
    'E18' = '  -1.00%  '
    'D19' = '17.08'
    'E19' = '  -1.00%  '
    'E20' = '  -0.13%  '
    'D21' = '29.268.38'
    'E21' = '  +0.18%  '
    'D23' = '11.10'
    'E23' = '  +0.14%  '
    'D24' = '2.098'
    'E24' = '  +1.76%  '
    'D25' = '156.41'
    'E25' = '  +0.00%  '
    'D26' = '19.42'
    'E26' = '  -0.72%  '
    'D27' = '5.973'
    'E27' = '  +1.85%  '
    'D28' = '118.64'
    'E28' = '  -0.86%  '
    'D29' = '1.910'
    'E29' = '  -4.99%  '
    'D30' = '0.09378'
    'E30' = '  -0.43%  '
    'D31' = '0.9054'
    'E31' = '  -2.17%  '
    'D32' = '5.290'
    'E32' = '  -1.53%  '
    'D33' = '1.328'
    'E33' = '  -1.38%  '
    'D34' = '3.217'
    'E34' = '  -1.61%  '
    'D35' = '1.186'
    'E35' = '  +1.79%  '
    'D36' = '0.05781'
    'E36' = '  -1.43%  '
    'D37' = '0.02092'
    'E37' = '  -0.41%  '
    'E38' = '  -0.26%  '
    'D39' = '7.750'
    'E39' = '  -3.71%  '
    'D40' = '0.5715'
    'E40' = '  -0.88%  '
    'D41' = '0.1788'
    'E41' = '  -1.43%  '
    'D42' = '9.769'
    'E42' = '  -2.65%  '
    'D43' = '11.96'
    'E43' = '  -1.17%  '
    'D44' = '0.5353'
    'E44' = '  -1.48%  '
    'D45' = '2.196'
    'E45' = '  -1.57%  '
    'D46' = '0.07042'
    'E46' = '  -0.83%  '
    'D47' = '1.853'
    'E47' = '  -1.53%  '
    'D48' = '2.572'
    'E48' = '  -0.11%  '
    'D49' = '113.37'
    'E50' = '  -3.62%  '
    'D51' = '71.24'
    'E51' = '  -0.79%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
